# Generate Report for Handback
# The handback transform for the "9cb99ef0-fa45-4deb-9f32-d7e4291f59fb" file
# failed for both target languages (zh-cn / de-de). Update the per-language
# status/error-detail columns accordingly and widen the "Error Detail" column
# so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 3 on each language sheet is the 9cb99ef0-... file.
# Status column (C) -> "Handback transform failed" (shared text, so both
# sheets pick it up since it's the very same value).
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# The Overview sheet mirrors each language's Status in columns E (zh-cn) and
# F (de-de) for this same row, so reflect the new status there too.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# Error Detail column (P) -> per-language transform error message.
$zhcn.Range("P3").Value = "Handback file name: yh5qzm1q.1in is different with handoff file name: 9cb99ef0-fa45-4deb-9f32-d7e4291f59fb.878cd327150e146efa3d4d885644520fb20e3f4d.zh-cn."
$dede.Range("P3").Value = "Handback file name: yh5qzm1q.1in is different with handoff file name: 9cb99ef0-fa45-4deb-9f32-d7e4291f59fb.878cd327150e146efa3d4d885644520fb20e3f4d.de-de."

# Widen the Error Detail column (P) on both sheets so the message is visible.
# (ColumnWidth is in "characters"; the engine persists OOXML <col width> as
# ColumnWidth + 5/6, matching Excel's own padding convention, so back that
# off here to land on an exact width="40" in the saved file.)
$finalColWidth = 40 - (5/6)
$zhcn.Columns.Item(16).ColumnWidth = $finalColWidth
$dede.Columns.Item(16).ColumnWidth = $finalColWidth
